# Updates cryptos list price/volume columns (D, E) to match refreshed market data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.383.17"
$ws.Range("E2").Value = "  +4.13%  "

$ws.Range("D3").Value = "3.496.03"
$ws.Range("E3").Value = "  +3.76%  "

$ws.Range("D5").Value = "'585.55"
$ws.Range("E5").Value = "  +2.83%  "

$ws.Range("D6").Value = "'147.46"
$ws.Range("E6").Value = "  +6.18%  "

$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("E8").Value = "  +1.31%  "

$ws.Range("D9").Value = "'7.71"
$ws.Range("E9").Value = "  +0.53%  "

$ws.Range("E10").Value = "  +4.25%  "

$ws.Range("D11").Value = "'0.399"
$ws.Range("E11").Value = "  +4.62%  "

$ws.Range("D12").Value = "4.097.56"
$ws.Range("E12").Value = "  +3.89%  "

$ws.Range("D13").Value = "'29.88"
$ws.Range("E13").Value = "  +7.74%  "

$ws.Range("E14").Value = "  -0.41%  "

$ws.Range("D15").Value = "3.503.89"
$ws.Range("E15").Value = "  +4.00%  "

$ws.Range("D16").Value = "'0.0000174"
$ws.Range("E16").Value = "  +4.13%  "

$ws.Range("D17").Value = "63.430.72"
$ws.Range("E17").Value = "  +4.07%  "

$ws.Range("E18").Value = "  +3.21%  "

$ws.Range("D19").Value = "'14.27"
$ws.Range("E19").Value = "  +5.36%  "

$ws.Range("D20").Value = "'9.49"
$ws.Range("E20").Value = "  +6.85%  "

$ws.Range("D21").Value = "'394.30"
$ws.Range("E21").Value = "  +3.25%  "

$ws.Range("D22").Value = "'0.565"
$ws.Range("E22").Value = "  +3.13%  "

$ws.Range("D23").Value = "'75.41"
$ws.Range("E23").Value = "  +0.00%  "

$ws.Range("E24").Value = "  -0.01%  "

$ws.Range("D25").Value = "'0.0000120"
$ws.Range("E25").Value = "  +8.89%  "

$ws.Range("D26").Value = "3.645.79"
$ws.Range("E26").Value = "  +4.03%  "

$ws.Range("E27").Value = "  -0.82%  "

$ws.Range("D28").Value = "'7.80"
$ws.Range("E28").Value = "  +9.06%  "

$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.10%  "

$ws.Range("D30").Value = "'8.25"
$ws.Range("E30").Value = "  +5.79%  "

$ws.Range("E31").Value = "  +2.69%  "

$ws.Range("E32").Value = "  +6.10%  "

$ws.Range("E33").Value = "  +0.01%  "

$ws.Range("D34").Value = "'23.83"
$ws.Range("E34").Value = "  +3.91%  "

$ws.Range("D35").Value = "'32.65"
$ws.Range("E35").Value = "  +29.49%  "

$ws.Range("E36").Value = "  +5.02%  "

$ws.Range("D37").Value = "'5.34"
$ws.Range("E37").Value = "  +8.60%  "

$ws.Range("D38").Value = "'172.16"
$ws.Range("E38").Value = "  +3.00%  "

$ws.Range("E39").Value = "  +9.15%  "

$ws.Range("D40").Value = "3.534.42"
$ws.Range("E40").Value = "  +3.80%  "

$ws.Range("D41").Value = "'0.0770"
$ws.Range("E41").Value = "  +1.18%  "

$ws.Range("D42").Value = "'0.803"
$ws.Range("E42").Value = "  +3.96%  "

$ws.Range("E43").Value = "  +7.78%  "

$ws.Range("D44").Value = "'4.51"
$ws.Range("E44").Value = "  +4.38%  "

$ws.Range("D45").Value = "'42.56"
$ws.Range("E45").Value = "  +0.26%  "

$ws.Range("E46").Value = "  +10.03%  "

$ws.Range("D47").Value = "2.611.36"
$ws.Range("E47").Value = "  +6.53%  "

$ws.Range("D48").Value = "'23.80"
$ws.Range("E48").Value = "  +7.41%  "

$ws.Range("D49").Value = "'2.26"
$ws.Range("E49").Value = "  +12.26%  "

$ws.Range("E50").Value = "  +2.56%  "

$ws.Range("D51").Value = "'0.0270"
$ws.Range("E51").Value = "  +5.11%  "
